$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Connector part numbers (designators filled first) ---
$ws.Range("A17").Value = "A30603-ND"
$ws.Range("B17").Value = "J1"
$ws.Range("B18").Value = "J2"
$ws.Range("A18").Value = "A30604-ND"
$ws.Range("A19").Value = "A30591-ND"
$ws.Range("A20").Value = "A30592-ND"

# --- Descriptions for the female/male connectors ---
$ws.Range("E17").Value = "4 pin female connector VAL-U-LOK"
$ws.Range("E18").Value = "6 pin female connector VAL-U-LOK"
$ws.Range("E19").Value = "4 pin male connector VAL-U-LOK"
$ws.Range("E20").Value = "6 pin male connector VAL-U-LOK"

# --- Crimp pin line ---
$ws.Range("A21").Value = "A99267CT-ND"
$ws.Range("E21").Value = "VAL-U-LOK crimp pin"

# --- Heatsink / thermal paste, moved down and re-described ---
$ws.Range("E26").Value = "Heatsink, optional"
$ws.Range("E27").Value = "Thermal paste, for heatsink"
$ws.Range("A27").Value = "345-1006-ND"

# --- Quantities, unit prices and line-total formulas ---
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 0.68
$ws.Range("I17").Formula = "=D17*C17"

$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 0.65
$ws.Range("I18").Formula = "=D18*C18"

$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0.21
$ws.Range("I19").Formula = "=D19*C19"

$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0.24
$ws.Range("I20").Formula = "=D20*C20"

$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 0.05
$ws.Range("I21").Formula = "=D21*C21"

# --- Heatsink row (A26 already held this part number - kept identical) ---
$ws.Range("A26").Value = "345-1092-ND‎ "
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1.46

# --- Thermal paste row ---
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 3.82

# --- Update the running total formula to cover the newly added rows ---
$ws.Range("K2").Formula = "=SUM(I2:I106)"

# --- Update the view: scroll down a bit and select A26, matching the edited document ---
$ws.Range("A26").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
